# Update gh-pages to output generated at 456a3b4
# F2 (想去人数) goes 32 -> 33, F3 (想去人数) goes 73 -> 76
# on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 33
    $ws.Range("F3").Value = 76
}
